$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4303
$ws.Range("L3").Value = 4566
$ws.Range("B4").Value = 1714
$ws.Range("I4").Value = 1841
$ws.Range("K4").Value = 1779
$ws.Range("L4").Value = 1132
$ws.Range("K6").Value = 9116
$ws.Range("L6").Value = 3933
$ws.Range("B7").Value = 23346
$ws.Range("I7").Value = 26311
$ws.Range("K7").Value = 27571
$ws.Range("L7").Value = 14194

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 273
$ws.Range("L3").Value = 313
$ws.Range("L4").Value = 68
$ws.Range("L6").Value = 255
$ws.Range("L7").Value = 942

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 101
$ws.Range("L3").Value = 123
$ws.Range("L7").Value = 315

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 184
$ws.Range("L4").Value = 38
$ws.Range("L7").Value = 658

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L6").Value = 144
$ws.Range("L7").Value = 522

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 238

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 115
$ws.Range("L7").Value = 472
$ws.Range("L8").Value = 942
$ws.Range("L10").Value = 95
$ws.Range("L15").Value = 108
$ws.Range("L19").Value = 395
$ws.Range("L20").Value = 359
$ws.Range("I27").Value = 225
$ws.Range("K27").Value = 260
$ws.Range("L27").Value = 127
$ws.Range("L29").Value = 795
$ws.Range("L33").Value = 658
$ws.Range("L37").Value = 522
$ws.Range("L46").Value = 35
$ws.Range("L47").Value = 103
$ws.Range("L51").Value = 175
$ws.Range("L53").Value = 167
$ws.Range("K54").Value = 525
$ws.Range("L60").Value = 86
$ws.Range("B63").Value = 418
$ws.Range("L67").Value = 485
$ws.Range("L76").Value = 214
$ws.Range("L77").Value = 96
$ws.Range("L78").Value = 184
$ws.Range("L79").Value = 374
$ws.Range("L83").Value = 315
$ws.Range("L85").Value = 738
$ws.Range("L86").Value = 111
$ws.Range("L89").Value = 203
$ws.Range("L90").Value = 143
$ws.Range("L91").Value = 196
$ws.Range("L92").Value = 42
$ws.Range("L94").Value = 178
$ws.Range("L96").Value = 155
$ws.Range("L99").Value = 238
$ws.Range("B101").Value = 23346
$ws.Range("I101").Value = 26311
$ws.Range("K101").Value = 27571
$ws.Range("L101").Value = 14194

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 187
$ws.Range("L4").Value = 36
$ws.Range("L7").Value = 485

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 275
$ws.Range("K7").Value = 525

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 241
$ws.Range("L6").Value = 205
$ws.Range("L7").Value = 795

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 123
$ws.Range("L7").Value = 395

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 43
$ws.Range("L3").Value = 38
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 83
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 126
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 374

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 115
$ws.Range("L7").Value = 359

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 158
$ws.Range("L3").Value = 150
$ws.Range("L6").Value = 121
$ws.Range("L7").Value = 472

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 72
$ws.Range("L7").Value = 178

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 57
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 40
$ws.Range("I4").Value = 31
$ws.Range("K4").Value = 32
$ws.Range("I7").Value = 225
$ws.Range("K7").Value = 260
$ws.Range("L7").Value = 127

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L4").Value = 12
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 143

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 53
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 28
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 298
$ws.Range("L7").Value = 738

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 96
